$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update O2 justification text
$ws.Range("O2").Value = "The report demonstrates strong evidence extraction quality with accurate citations and full sentence quotations, though it lacks some depth in source diversity. Coverage of proxy and representation dimensions is comprehensive, detailing definitions, examples, and methods, but could benefit from more diverse examples. Structure and formatting are clear and consistent, aiding readability. The relevance and faithfulness of the content are maintained, with no unsupported assumptions. However, the identification of missing disclosures is somewhat limited, lacking specific details on OpenAI's practices. Audit usefulness is reasonable, but the report could improve by providing more actionable insights and clearer audit trails. Overall, the report is strong but has room for improvement in depth and specificity."

# Update row 3 numeric values
$ws.Range("B3").Value = 3
$ws.Range("C3").Value = 4
$ws.Range("D3").Value = 3
$ws.Range("H3").Value = 0.75
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 0.45
$ws.Range("N3").Value = 3.4

# Update O3 justification text
$ws.Range("O3").Value = "The report provides a reasonable extraction of evidence with accurate citations, but lacks full sentence quotations in some cases, leading to a score of 3 for evidence extraction quality. Coverage of proxy dimensions is strong, with clear definitions and examples, but lacks depth in methods, earning a 4. Structure is adequate but could be more organized, scoring a 3. The relevance is high, with evidence grounded in sources, scoring a 4. Missing disclosures are identified, but the analysis is not exhaustive, resulting in a 3. Overall audit usefulness is moderate due to some vagueness, scoring a 3. A key weakness is the lack of detailed methods for detecting proxy effects."
